# Removing TWB PNPC table and adding two extra fields to TWB Episode instead.

$wb = $excel.ActiveWorkbook

# 1. Delete the "TWB PNPCs" worksheet entirely.
$pnpcSheet = $wb.Worksheets.Item("TWB PNPCs")
$pnpcSheet.Delete()

# 2. Add two new fields to "TWB Episodes":
#    twb_primary_nominated_professional_contact_entry_date
#    twb_primary_nominated_professional_contact_exit_date
# These are inserted right after "twb_primary_nominated_professional_consent_date"
# (column K) and before "twb_previous_suicide_attempts" (old column L), so the
# old L/M columns (previous_suicide_attempts / method_of_suicide_attempt) shift
# right into new columns N/O.
$twbEpisodes = $wb.Worksheets.Item("TWB Episodes")

$twbEpisodes.Range("L1:M1").EntireColumn.Insert()

$twbEpisodes.Range("L1").Value = "twb_primary_nominated_professional_contact_entry_date"
$twbEpisodes.Range("M1").Value = "twb_primary_nominated_professional_contact_exit_date"

$twbEpisodes.Range("L2").Value = 16042020
$twbEpisodes.Range("M2").Value = 9099999

$twbEpisodes.Range("L3").Value = 9099999
$twbEpisodes.Range("M3").Value = 9099999
